$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = <new price or $null>; E = <new volume pct> }
$updates = @(
    @{ Row = 2; D = "27.537.10"; E = "  -0.30%  " }
    @{ Row = 3; D = "1.646.61"; E = "  -1.13%  " }
    @{ Row = 5; D = "212.58"; E = "  -1.37%  " }
    @{ Row = 6; D = "0.529"; E = "  +3.86%  " }
    @{ Row = 7; D = "0.999"; E = "  -0.07%  " }
    @{ Row = 8; D = "23.57"; E = "  -2.58%  " }
    @{ Row = 9; D = "0.258"; E = "  -1.96%  " }
    @{ Row = 10; D = "0.0612"; E = "  -1.37%  " }
    @{ Row = 11; D = "0.0891"; E = "  +1.53%  " }
    @{ Row = 12; D = "1.879.10"; E = "  -1.17%  " }
    @{ Row = 13; D = "1.643.51"; E = "  -1.72%  " }
    @{ Row = 14; D = "0.592"; E = "  +3.74%  " }
    @{ Row = 15; D = $null; E = "  -2.27%  " }
    @{ Row = 16; D = "64.55"; E = "  -2.64%  " }
    @{ Row = 17; D = "27.504.38"; E = "  -0.35%  " }
    @{ Row = 18; D = "232.05"; E = "  -4.13%  " }
    @{ Row = 19; D = $null; E = "  -1.03%  " }
    @{ Row = 20; D = "7.55"; E = "  -1.78%  " }
    @{ Row = 21; D = $null; E = "  -0.03%  " }
    @{ Row = 22; D = $null; E = "  -3.85%  " }
    @{ Row = 23; D = "9.76"; E = "  +4.14%  " }
    @{ Row = 24; D = $null; E = "  -1.35%  " }
    @{ Row = 25; D = "148.36"; E = "  +1.30%  " }
    @{ Row = 26; D = "7.05"; E = "  -3.09%  " }
    @{ Row = 27; D = "0.113"; E = "  +1.66%  " }
    @{ Row = 28; D = $null; E = "  -0.06%  " }
    @{ Row = 29; D = "15.64"; E = "  -4.54%  " }
    @{ Row = 30; D = $null; E = "  -2.87%  " }
    @{ Row = 31; D = $null; E = "  -3.32%  " }
    @{ Row = 32; D = $null; E = "  -0.91%  " }
    @{ Row = 33; D = "3.17"; E = "  +1.18%  " }
    @{ Row = 34; D = "1.426.75"; E = "  -2.14%  " }
    @{ Row = 35; D = $null; E = "  +0.22%  " }
    @{ Row = 36; D = "2.37"; E = "  +0.21%  " }
    @{ Row = 37; D = "0.569"; E = "  -1.22%  " }
    @{ Row = 38; D = $null; E = "  -4.41%  " }
    @{ Row = 39; D = $null; E = "  -3.36%  " }
    @{ Row = 40; D = $null; E = "  -0.82%  " }
    @{ Row = 41; D = "0.999"; E = "  -0.06%  " }
    @{ Row = 42; D = $null; E = "  +3.13%  " }
    @{ Row = 43; D = "5.56"; E = "  +2.52%  " }
    @{ Row = 44; D = $null; E = "  -1.94%  " }
    @{ Row = 45; D = "2.25"; E = "  +0.78%  " }
    @{ Row = 46; D = "65.12"; E = "  -6.89%  " }
    @{ Row = 47; D = "1.788.84"; E = "  -1.10%  " }
    @{ Row = 48; D = $null; E = "  -2.55%  " }
    @{ Row = 49; D = "88.23"; E = "  -0.77%  " }
    @{ Row = 50; D = "0.0₆0108"; E = "  -0.33%  " }
    @{ Row = 51; D = $null; E = "  -1.43%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $isNumericLooking = $u.D -match "^[0-9]*\.?[0-9]+$"
        if ($isNumericLooking) {
            # Force text storage so Excel does not coerce the digit-dot string to a number
            $cell.NumberFormat = "@"
            $cell.Value = $u.D
            $cell.ClearFormats()
        } else {
            $cell.Value = $u.D
        }
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
